# Update "Förändrad" date column (C) for rows 2-24 from 45203 to 45204.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Avverkningsanmälningar")

for ($row = 2; $row -le 24; $row++) {
    $cell = $ws.Cells.Item($row, 3)
    if ($cell.Value2 -eq 45203) {
        $cell.Value2 = 45204
    }
}
